$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 53
$ws.Range("H53").Value = 227.94737
$ws.Range("I53").Value = 149.61539
$ws.Range("K53").Value = 149.61539
$ws.Range("M53").Value = 487.38461
# Row 62
$ws.Range("H62").Value = 7782.375
$ws.Range("I62").Value = 6867
$ws.Range("K62").Value = 6867
$ws.Range("M62").Value = -6243
# Row 65
$ws.Range("H65").Value = 7782.375
$ws.Range("I65").Value = 6867
$ws.Range("K65").Value = 34335
$ws.Range("M65").Value = -31215
# Row 98
$ws.Range("H98").Value = 3143.2354
$ws.Range("I98").Value = 1179.6666
$ws.Range("K98").Value = 1179.6666
$ws.Range("M98").Value = 318.3334
# Row 116
$ws.Range("H116").Value = 5363.353
$ws.Range("I116").Value = 3653
$ws.Range("K116").Value = 3653
$ws.Range("M116").Value = -211
# Row 122
$ws.Range("H122").Value = 3143.2354
$ws.Range("I122").Value = 1179.6666
$ws.Range("K122").Value = 3538.9998
$ws.Range("M122").Value = -1088.9998
# Row 132
$ws.Range("H132").Value = 814.9048
$ws.Range("I132").Value = 791.6842
$ws.Range("K132").Value = 2375.0526
$ws.Range("M132").Value = 154.9474
# Row 137
$ws.Range("H137").Value = 1626.6428
$ws.Range("I137").Value = 1693.5834
$ws.Range("J137").Value = 1225
$ws.Range("K137").Value = 5080.7502
$ws.Range("L137").Value = 3675
$ws.Range("M137").Value = -2530.7502
$ws.Range("N137").Value = -8775

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 74
$ws.Range("H74").Value = 2098.5
$ws.Range("I74").Value = 2219.875
$ws.Range("J74").Value = 1936.6666
$ws.Range("K74").Value = 2219.875
$ws.Range("L74").Value = 1936.6666
$ws.Range("M74").Value = -1345.875
$ws.Range("N74").Value = -3684.6666
# Row 77
$ws.Range("H77").Value = 2098.5
$ws.Range("I77").Value = 2219.875
$ws.Range("J77").Value = 1936.6666
$ws.Range("K77").Value = 11099.375
$ws.Range("L77").Value = 9683.333000000001
$ws.Range("M77").Value = -6731.375
$ws.Range("N77").Value = -18419.333

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 22
$ws.Range("H22").Value = 612.55554
$ws.Range("I22").Value = 484.125
$ws.Range("K22").Value = 484.125
$ws.Range("M22").Value = -311.125
# Row 86
$ws.Range("H86").Value = 3158.8572
$ws.Range("I86").Value = 350.83334
$ws.Range("J86").Value = 20007
$ws.Range("K86").Value = 350.83334
$ws.Range("L86").Value = 20007
$ws.Range("M86").Value = 772.16666
$ws.Range("N86").Value = -22253
# Row 89
$ws.Range("H89").Value = 3158.8572
$ws.Range("I89").Value = 350.83334
$ws.Range("J89").Value = 20007
$ws.Range("K89").Value = 1754.1667
$ws.Range("L89").Value = 100035
$ws.Range("M89").Value = 3861.8333
$ws.Range("N89").Value = -111267
# Row 99
$ws.Range("H99").Value = 2069.25
$ws.Range("I99").Value = 2086.4546
$ws.Range("K99").Value = 2086.4546
$ws.Range("M99").Value = -588.4546

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 58
$ws.Range("H58").Value = 2199.3157
$ws.Range("I58").Value = 2164.2354
$ws.Range("J58").Value = 2497.5
$ws.Range("K58").Value = 2164.2354
$ws.Range("L58").Value = 2497.5
$ws.Range("M58").Value = -1961.2354
$ws.Range("N58").Value = -2903.5
# Row 104
$ws.Range("H104").Value = 44745
$ws.Range("J104").Value = 44745
$ws.Range("L104").Value = 44745
$ws.Range("N104").Value = -49987
# Row 122
$ws.Range("H122").Value = 2154.2
$ws.Range("I122").Value = 1654
$ws.Range("J122").Value = 2904.5
$ws.Range("K122").Value = 4962
$ws.Range("L122").Value = 8713.5
$ws.Range("M122").Value = -2512
$ws.Range("N122").Value = -13613.5
# Row 124
$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").ClearContents()
# Row 136
$ws.Range("H136").Value = 2199.3157
$ws.Range("I136").Value = 2164.2354
$ws.Range("J136").Value = 2497.5
$ws.Range("K136").Value = 6492.706200000001
$ws.Range("L136").Value = 7492.5
$ws.Range("M136").Value = -3942.706200000001
$ws.Range("N136").Value = -12592.5

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 4
$ws.Range("H4").Value = 7692411.5
$ws.Range("I4").Value = 7692411.5
$ws.Range("K4").Value = 23077234.5
$ws.Range("M4").Value = -23077122.5
# Row 37
$ws.Range("H37").Value = 79977.5
$ws.Range("J37").Value = 79977.5
$ws.Range("L37").Value = 239932.5
$ws.Range("N37").Value = -240156.5
# Row 130
$ws.Range("H130").Value = 1915.1666
$ws.Range("I130").Value = 1936.6666
$ws.Range("J130").Value = 1893.6666
$ws.Range("K130").Value = 5809.9998
$ws.Range("L130").Value = 5680.9998
$ws.Range("M130").Value = -789.9997999999996
$ws.Range("N130").Value = -15720.9998

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 7768.615
$ws.Range("I70").Value = 7790.636
$ws.Range("K70").Value = 7790.636
$ws.Range("M70").Value = -7520.636
# Row 73
$ws.Range("H73").Value = 7768.615
$ws.Range("I73").Value = 7790.636
$ws.Range("K73").Value = 7790.636
$ws.Range("M73").Value = -6854.636
# Row 135
$ws.Range("H135").Value = 285657.5
$ws.Range("J135").Value = 285657.5
$ws.Range("L135").Value = 285657.5
$ws.Range("N135").Value = -295797.5

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 61
$ws.Range("H61").Value = 5967.2
$ws.Range("I61").Value = 6459.25
$ws.Range("K61").Value = 6459.25
$ws.Range("M61").Value = -6257.25
# Row 113
$ws.Range("H113").Value = 5967.2
$ws.Range("I113").Value = 6459.25
$ws.Range("K113").Value = 6459.25
$ws.Range("M113").Value = -4289.25
# Row 132
$ws.Range("H132").Value = 3276.889
$ws.Range("I132").Value = 2531.1667
$ws.Range("J132").Value = 4768.3335
$ws.Range("K132").Value = 7593.500100000001
$ws.Range("L132").Value = 14305.0005
$ws.Range("M132").Value = -5063.500100000001
$ws.Range("N132").Value = -19365.0005
# Row 136
$ws.Range("H136").Value = 3068.4119
$ws.Range("I136").Value = 2468.0908
$ws.Range("K136").Value = 7404.2724
$ws.Range("M136").Value = -4854.2724

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 109
$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()
# Row 122
$ws.Range("H122").Value = 9254
$ws.Range("J122").Value = 11318
$ws.Range("L122").Value = 33954
$ws.Range("N122").Value = -38854
# Row 126
$ws.Range("H126").Value = 5567.6
$ws.Range("I126").Value = 4113.1113
$ws.Range("J126").Value = 7749.3335
$ws.Range("K126").Value = 12339.3339
$ws.Range("L126").Value = 23248.0005
$ws.Range("M126").Value = -9869.333899999998
$ws.Range("N126").Value = -28188.0005

Write-Host "Applied scheduled market data update to 8 sheets (32 rows)."
